$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style of the existing
# header row (bold, centered, bordered - same as H1's style).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

$i_vals = @(3,7,4,5,6,6,2,7,2,1,6,5,7,7,5,9,4,6,5,2,6,8,6,6,10,6,8,7,8,6,8,3,8,6,5,6,5)
$j_vals = @(5,8,6,6,8,6,4,8,3,3,6,7,8,8,6,9,6,7,6,4,6,8,8,8,10,8,10,8,9,8,8,5,9,7,6,7,6)

for ($k = 0; $k -lt $i_vals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $i_vals[$k]
    $ws.Cells.Item($row, 10).Value = $j_vals[$k]
}
